# Updates the crypto price/volume table on Sheet1 (rows 2-51) to reflect the
# latest scrape. Price cells in column D are plain decimal-looking strings
# (e.g. "1.000", "18.00") that Excel would otherwise silently coerce to
# numbers (dropping the meaningful trailing zeros / thousands-style dots),
# so each one is written with a temporary Text number format and then reset
# back to the "Normal" style to avoid leaving stray formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.094.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.834.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6799"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07451"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.018"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6757"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.130.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008254"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.075.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.338"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1441"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.701"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.246"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.126"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05404"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7528"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.307.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01818"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.717"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9330"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.072"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08583"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +34.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.981.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5180"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.432"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.768"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000120"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.97%  "
